$wb = $excel.ActiveWorkbook

# Sheet "存款" (deposits)
$ws1 = $wb.Worksheets.Item("存款")
$ws1.Range("B2").Value = "中國信託商業銀行斗六分行"
$ws1.Range("F5").Value = "3641580"
$ws1.Range("F9").Value = "2506288"

# Sheet "其他有價證券" (other securities / property)
$ws2 = $wb.Worksheets.Item("其他有價證券")
$ws2.Range("B2").Value = "(九）珠寶古董字畫及#"
$ws2.Range("C2").Value = "他具有相當價值之財產（總作"
$ws2.Range("E2").Value = "1額：新臺幣元）"
$ws2.Range("B3").Value = "財產種類"
$ws2.Range("C3").Value = "項"
